# Update odds for the Argentina match (row 2) that stayed in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.5
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("AI2").Value = 17
$ws.Range("AJ2").Value = 13
$ws.Range("AK2").Value = 41
$ws.Range("AN2").Value = 4
$ws.Range("AX2").Value = 21
$ws.Range("BD2").Value = 151

# The Bolivia match (old row 3) was removed entirely; this shifts the
# Colombia match up from row 4 to row 3, and the Paraguay match up from
# row 5 to row 4.
$ws.Rows(3).Delete()

# The Paraguay match (now row 4) was removed entirely too.
$ws.Rows(4).Delete()

# A handful of odds on the Colombia match (now row 3) were also refreshed.
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.73
$ws.Range("X3").Value = 8
$ws.Range("AC3").Value = 8
$ws.Range("AM3").Value = 41
$ws.Range("AN3").Value = 3.75
$ws.Range("AR3").Value = 51
$ws.Range("AY3").Value = 34
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 301
